$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-13 from 45204 (2023-10-05)
# to 45207 (2023-10-08), keeping the existing date number format/style.
$ws.Range("C2:C13").Value = 45207
